$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the numeric-looking Price cells as Text so Excel does not
# auto-convert them to numbers (which would also strip trailing zeros).
$ws.Range("D4,D5,D6,D7,D8,D9,D10,D11,D13,D14,D15,D16,D18,D19,D20,D22,D23,D24,D25,D26,D27,D28,D29,D30,D31,D32,D33,D34,D35,D36,D37,D38,D39,D40,D41,D42,D43,D44,D45,D46,D47,D48,D49,D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.405.71"
$ws.Range("E2").Value = "  -0.96%  "
$ws.Range("D3").Value = "1.916.79"
$ws.Range("E3").Value = "  +2.02%  "
$ws.Range("D4").Value = "0.9996"
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").Value = "240.99"
$ws.Range("E5").Value = "  +1.69%  "
$ws.Range("D6").Value = "0.9996"
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("D7").Value = "0.4672"
$ws.Range("E7").Value = "  -1.78%  "
$ws.Range("D8").Value = "0.2850"
$ws.Range("E8").Value = "  +0.97%  "
$ws.Range("D9").Value = "0.06814"
$ws.Range("E9").Value = "  +5.09%  "
$ws.Range("D10").Value = "107.29"
$ws.Range("E10").Value = "  +12.70%  "
$ws.Range("D11").Value = "18.12"
$ws.Range("E11").Value = "  -2.74%  "
$ws.Range("D12").Value = "1.909.86"
$ws.Range("E12").Value = "  +1.25%  "
$ws.Range("D13").Value = "0.07625"
$ws.Range("E13").Value = "  +0.94%  "
$ws.Range("D14").Value = "5.175"
$ws.Range("E14").Value = "  +2.34%  "
$ws.Range("D15").Value = "0.6549"
$ws.Range("E15").Value = "  +1.04%  "
$ws.Range("D16").Value = "287.38"
$ws.Range("E16").Value = "  -3.45%  "
$ws.Range("D17").Value = "30.421.19"
$ws.Range("E17").Value = "  -0.96%  "
$ws.Range("B18").Value = "Avalanche"
$ws.Range("C18").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D18").Value = "12.99"
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("D19").Value = "0.000007607"
$ws.Range("E19").Value = "  +2.06%  "
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").Value = "0.9993"
$ws.Range("E20").Value = "  -0.27%  "
$ws.Range("D21").Value = "2.159.03"
$ws.Range("E21").Value = "  +0.94%  "
$ws.Range("D22").Value = "0.9997"
$ws.Range("E22").Value = "  -0.35%  "
$ws.Range("D23").Value = "5.213"
$ws.Range("E23").Value = "  +1.88%  "
$ws.Range("D24").Value = "6.194"
$ws.Range("E24").Value = "  +1.53%  "
$ws.Range("D25").Value = "168.02"
$ws.Range("E25").Value = "  -0.64%  "
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "21.53"
$ws.Range("E26").Value = "  +10.75%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "9.263"
$ws.Range("E27").Value = "  +0.46%  "
$ws.Range("D28").Value = "2.049"
$ws.Range("E28").Value = "  +5.14%  "
$ws.Range("D29").Value = "0.1070"
$ws.Range("E29").Value = "  +1.31%  "
$ws.Range("D30").Value = "1.374"
$ws.Range("E30").Value = "  +1.40%  "
$ws.Range("D31").Value = "4.133"
$ws.Range("E31").Value = "  -0.37%  "
$ws.Range("D32").Value = "3.940"
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("D33").Value = "0.05038"
$ws.Range("E33").Value = "  +1.15%  "
$ws.Range("D34").Value = "0.7372"
$ws.Range("E34").Value = "  +2.77%  "
$ws.Range("D35").Value = "1.147"
$ws.Range("E35").Value = "  -1.63%  "
$ws.Range("D36").Value = "0.9987"
$ws.Range("E36").Value = "  -0.19%  "
$ws.Range("D37").Value = "2.732"
$ws.Range("E37").Value = "  +0.66%  "
$ws.Range("D38").Value = "0.02030"
$ws.Range("E38").Value = "  +5.75%  "
$ws.Range("D39").Value = "2.687"
$ws.Range("E39").Value = "  -0.46%  "
$ws.Range("D40").Value = "2.051"
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("D41").Value = "108.87"
$ws.Range("E41").Value = "  +2.06%  "
$ws.Range("D42").Value = "0.8731"
$ws.Range("E42").Value = "  -2.40%  "
$ws.Range("D43").Value = "5.824"
$ws.Range("E43").Value = "  +4.71%  "
$ws.Range("D44").Value = "52.95"
$ws.Range("E44").Value = "  +27.12%  "
$ws.Range("D45").Value = "0.9993"
$ws.Range("E45").Value = "  -0.19%  "
$ws.Range("D46").Value = "0.4199"
$ws.Range("E46").Value = "  +0.82%  "
$ws.Range("D47").Value = "67.56"
$ws.Range("E47").Value = "  +2.91%  "
$ws.Range("D48").Value = "7.149"
$ws.Range("E48").Value = "  -2.08%  "
$ws.Range("D49").Value = "9.204"
$ws.Range("E49").Value = "  +4.02%  "
$ws.Range("E50").Value = "  -0.37%  "
$ws.Range("D51").Value = "34.65"
$ws.Range("E51").Value = "  +0.79%  "
